# Apply the "automatic update of files" commit: the five data rows (2-6)
# got re-ordered/re-fetched from source, so each row's content fields are
# replaced with the values belonging to a different underlying record.
# Row-to-row content mapping (target row <- source-of-truth row, i.e. what
# currently sits in that row before the edit):
#   2 <- 4
#   3 <- 5
#   4 <- 3
#   5 <- 6
#   6 <- 2

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- capture current ("before") values for the columns that rotate ------
# NOTE: use .Value2 (not .Value) when round-tripping values through
# variables/hashtables - .Value can surface an unresolved property
# descriptor instead of the real scalar in this host.
$cols = @("A","B","D","E","F","G","H","Q","R")
$before = @{}
foreach ($r in 2..6) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $before[$r] = $rowVals
}

$rowMap = @{ 2 = 4; 3 = 5; 4 = 3; 5 = 6; 6 = 2 }

foreach ($target in 2..6) {
    $source = $rowMap[$target]
    $src = $before[$source]
    foreach ($c in $cols) {
        $ws.Range("$c$target").Value2 = $src[$c]
    }
}

# --- move the "blomning"/"Blommande" activity block from row 2 to row 6 -
# Row 2 currently carries J2 (blank), K2="blomning", L2 (blank), N2 (blank),
# AC2="Blommande", AF2 (blank). After the edit these belong to row 6, and
# row 2 no longer has any of these cells.

$ws.Range("K6").Value2 = $ws.Range("K2").Value2
$ws.Range("AC6").Value2 = $ws.Range("AC2").Value2

# re-create the blank placeholder cells on row 6 (J6, L6, N6, AF6) by
# copying an already-blank inline-string cell (I2) so the cell exists
# without altering style/number-format.
$ws.Range("I2").Copy($ws.Range("J6"))
$ws.Range("I2").Copy($ws.Range("L6"))
$ws.Range("I2").Copy($ws.Range("N6"))
$ws.Range("I2").Copy($ws.Range("AF6"))

# remove the now-obsolete cells from row 2
$ws.Range("J2").ClearContents()
$ws.Range("K2").ClearContents()
$ws.Range("L2").ClearContents()
$ws.Range("N2").ClearContents()
$ws.Range("AC2").ClearContents()
$ws.Range("AF2").ClearContents()
